$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows at 4 and 5 for new simulated items "Holden" and "Rizzie Spiral" ---
# This pushes the old rows 4..29 (and their C:T simulation payloads) down to rows 6..31.
$ws.Range("4:5").Insert()

# Copy formatting (bold, border, centered) from row 3 col A onto the new A4:A5 cells
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 4: new item "Holden" ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.216722602658349
$ws.Range("D4").Value = 0.7965150409393728
$ws.Range("E4").Value = 0.7965150409393728
$ws.Range("F4").Value = 0.8445237732028114
$ws.Range("G4").Value = 0.8445237732028114
$ws.Range("H4").Value = 0.8276502630463403
$ws.Range("I4").Value = 1.404298928318148
$ws.Range("J4").Value = 1.004397384289178
$ws.Range("K4").Value = 0.8445237732028114
$ws.Range("L4").Value = 1.216722602658349
$ws.Range("M4").Value = 1.006618821798861
$ws.Range("N4").Value = 1.006618821798861
$ws.Range("O4").Value = 0.9469626355480205
$ws.Range("P4").Value = 0.9525871389335109
$ws.Range("Q4").Value = 0.9525871389335109
$ws.Range("R4").Value = 0.925571297500836
$ws.Range("S4").Value = 0.925571297500836
$ws.Range("T4").Value = 1.015684665409033

# --- Row 5: new item "Rizzie Spiral" ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.266896308884938
$ws.Range("D5").Value = 0.0008677191756345288
$ws.Range("E5").Value = 0.0008677191756345288
$ws.Range("F5").Value = 2.995904087481815
$ws.Range("G5").Value = 2.995904087481815
$ws.Range("H5").Value = 0.00218331654402714
$ws.Range("I5").Value = 0.03040335315894161
$ws.Range("J5").Value = 1.410864304955304
$ws.Range("K5").Value = 2.995904087481815
$ws.Range("L5").Value = 1.266896308884938
$ws.Range("M5").Value = 0.6338820140302864
$ws.Range("N5").Value = 0.6338820140302864
$ws.Range("O5").Value = 0.4233157815348667
$ws.Range("P5").Value = 1.421222705180796
$ws.Range("Q5").Value = 1.421222705180796
$ws.Range("R5").Value = 1.814893050756051
$ws.Range("S5").Value = 1.814893050756051
$ws.Range("T5").Value = 0.9511865150334436

# --- The rows pushed down keep their OLD column-A index; renumber it to match the new row ---
For ($r = 6; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Rename "Thomas Hex" to "Matthies Hex" (now at row 11 after the insert) ---
$ws.Range("B11").Value = "Matthies Hex"

Write-Host "done"